$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.52879536151886
$ws.Range("B1").Value = 3.822039365768433
$ws.Range("C1").Value = 3.016901731491089
$ws.Range("D1").Value = 2.4290452003479
$ws.Range("E1").Value = 1.420446038246155
